$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ------------------------------------------------------------------
# Row 37: new "Installation error ..." bug, logged by AD on 31 Jan
# 2019 (serial 43496), resolved 1 Feb 2019 (serial 43497). Column C
# carries rich text: the main description in the default font, with
# the resolution note appended in red. The date columns reuse the
# existing date format (numFmtId 14) by copying it from an existing
# date cell instead of typing a NumberFormat string, so no redundant
# number format gets created.
# ------------------------------------------------------------------
$ws.Range("A37").Value = "AD"

$ws.Range("B37").Value = 43496
$ws.Range("B2").Copy()
$ws.Range("B37").PasteSpecial(-4122)

$part1 = 'Installation error, "index.html": No such file or directory. '
$part2 = "Solved by deleting build folder in tar package."
$ws.Range("C37").Value = ($part1 + $part2)

# Touch the whole cell's font red once so the red font is registered
# in the workbook's shared font table, then put the cell style back
# to Normal - only the second run (set via Characters below) should
# actually render red.
$ws.Range("C37").Font.Color = 255
$ws.Range("C37").Style = "Normal"

$part1Len = $part1.Length
$part2Len = $part2.Length
$redRun = $ws.Range("C37").Characters($part1Len + 1, $part2Len)
$redRun.Font.Color = 255

$ws.Range("D37").Value = 43497
$ws.Range("D3").Copy()
$ws.Range("D37").PasteSpecial(-4122)

# ------------------------------------------------------------------
# Row 38: second new bug about sumby + ordered factors, raised 1 Feb
# 2019 and resolved the same day.
# ------------------------------------------------------------------
$ws.Range("A38").Value = "AD"

$ws.Range("B38").Value = 43497
$ws.Range("B2").Copy()
$ws.Range("B38").PasteSpecial(-4122)

$ws.Range("C38").Value = "Sumby doesn't work for ordered factor. "

$ws.Range("D38").Value = 43497
$ws.Range("D3").Copy()
$ws.Range("D38").PasteSpecial(-4122)

# ------------------------------------------------------------------
# Scroll / selection bookkeeping to match where Excel would leave the
# cursor after typing the two new rows in.
# ------------------------------------------------------------------
$excel.ActiveWindow.ScrollRow = 13
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("C39").Select()
